$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells D1/E1: drop the red "*" rich-text run, keep plain "Gender"/"ROLE" ---
$ws.Range("D1").Value2 = "Gender"
$ws.Range("E1").Value2 = "ROLE"

# --- Clear the sample data rows 3 and 4 (A:E), keeping their existing cell styles ---
$ws.Range("A3:E4").ClearContents()

# --- Hyperlinks: only A2's mailto link should remain; drop the ones on A3/A4.
#     This engine only supports deleting the whole Hyperlinks collection, so
#     delete everything then re-add the single hyperlink that should survive. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:ANV@gmail.com")

# --- New column L width (bestFit-style width added to the <cols> definition) ---
$ws.Columns.Item(12).ColumnWidth = 18

# --- Selection moves to F14 ---
$ws.Range("F14").Select()
